# Features_Overview.xlsx — "4.4 Service Coordination" -> "4.4 Service Workflows"
# The row for this characteristic (row 16) previously only had data in column H;
# the rest of the row is now filled in with support info for every framework.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the characteristic itself.
$ws.Cells.Item(16, 1).Value = "4.4 Service Workflows"

# Columns B (Axon) and C (Eventuate) are now marked as supported ("✓").
# Re-use the formatting already used for a "✓" cell elsewhere in the table
# (D4) so the new cells pick up the same font/fill as the rest of the sheet.
$checkSrc = $ws.Cells.Item(4, 4)
$checkSrc.Copy($ws.Cells.Item(16, 2))
$checkSrc.Copy($ws.Cells.Item(16, 3))

# Columns D through O (Helidon, Lagom, Micronaut, Spring Cloud, MicroProfile,
# Falcon, Nameko, Kit, Micro, NODE.JS/Devis, Moleculer) are now marked as not
# supported ("✗"), matching the formatting already used for a "✗" cell (I5).
$crossSrc = $ws.Cells.Item(5, 9)
$crossSrc.Copy($ws.Cells.Item(16, 4))
$crossSrc.Copy($ws.Cells.Item(16, 5))
$crossSrc.Copy($ws.Cells.Item(16, 6))
$crossSrc.Copy($ws.Cells.Item(16, 7))
$crossSrc.Copy($ws.Cells.Item(16, 9))
$crossSrc.Copy($ws.Cells.Item(16, 10))
$crossSrc.Copy($ws.Cells.Item(16, 11))
$crossSrc.Copy($ws.Cells.Item(16, 12))
$crossSrc.Copy($ws.Cells.Item(16, 13))
$crossSrc.Copy($ws.Cells.Item(16, 14))
$crossSrc.Copy($ws.Cells.Item(16, 15))

# Leave the cursor where the author ended up after editing the row.
$ws.Range("Q16").Select()
